$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LogInPage")
$ws.Range("H89").NumberFormat = "mm-dd-yy"
$ws.Range("H89").Formula = "'12/12/1987"
$v = $ws.Range("H89").Value2
Write-Host "H89 value2:" $v
